$wb = $excel.ActiveWorkbook

# --- Schedule sheet updates ---
$ws = $wb.Worksheets.Item("Schedule")
$ws.Range("E3").Value = 367.9627575000001
$ws.Range("F3").Value = 24.33616121031747
$ws.Range("A4").Value = 46039.29166666666
$ws.Range("B4").Value = 46039.45833333334
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 15.12
$ws.Range("E4").Value = 177.4800885
$ws.Range("F4").Value = 11.73810109126984
$ws.Range("A5").Value = 46039.5625
$ws.Range("B5").Value = 46039.89583333334
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 30.24
$ws.Range("E5").Value = 74.748531
$ws.Range("F5").Value = 2.471842956349207

# --- Detailed sheet updates ---
$ws2 = $wb.Worksheets.Item("Detailed")
$ws2.Range("B43").Value = 29.85322
$ws2.Range("B44").Value = 0.85459
$ws2.Range("B45").Value = 65
$ws2.Range("C45").Value = "historical"
$ws2.Range("B46").Value = 65
$ws2.Range("C46").Value = "historical"
$ws2.Range("B47").Value = 64.8901
$ws2.Range("B48").Value = 64.8901
$ws2.Range("B49").Value = 64.8901
$ws2.Range("B50").Value = 57.06003
$ws2.Range("B51").Value = 57.06003
$ws2.Range("B52").Value = 57.06003
$ws2.Range("B53").Value = 57.06003
$ws2.Range("B54").Value = 40.54
$ws2.Range("B57").Value = 36.06
$ws2.Range("B58").Value = 36.06
$ws2.Range("B59").Value = 57.3
$ws2.Range("B60").Value = 57.06017
$ws2.Range("B61").Value = 57.3
$ws2.Range("B63").Value = 36.2
$ws2.Range("E64").Value = "ON"
$ws2.Range("B65").Value = 0.7
$ws2.Range("E65").Value = "ON"
$ws2.Range("B66").Value = 0.51003
$ws2.Range("B67").Value = 0.7
$ws2.Range("B68").Value = 36.06011
$ws2.Range("B69").Value = 35.88
$ws2.Range("B70").Value = 36.06011
$ws2.Range("B71").Value = 36.06032
$ws2.Range("B72").Value = 36.06046
$ws2.Range("E72").Value = "OFF"
$ws2.Range("B73").Value = 36.06046
$ws2.Range("E73").Value = "OFF"
$ws2.Range("B74").Value = 28.43746
$ws2.Range("E74").Value = "OFF"
$ws2.Range("B75").Value = 27.73
$ws2.Range("B76").Value = 36.06045
$ws2.Range("B77").Value = 2.47275
$ws2.Range("E77").Value = "ON"
$ws2.Range("B78").Value = 0.7
$ws2.Range("E78").Value = "ON"
$ws2.Range("B79").Value = 7.62095
$ws2.Range("B80").Value = 5.3077
$ws2.Range("B81").Value = 15.48773
$ws2.Range("B82").Value = 30.34528
$ws2.Range("B83").Value = 0.3475
$ws2.Range("B85").Value = -3.03541
$ws2.Range("B86").Value = -6.20557
$ws2.Range("B87").Value = -3.30783
$ws2.Range("B88").Value = 0.01089
$ws2.Range("B89").Value = 2.23907
$ws2.Range("B91").Value = 8.61159
$ws2.Range("B92").Value = 8.47514
$ws2.Range("B93").Value = 4.42107
$ws2.Range("E93").Value = "OFF"
$ws2.Range("B94").Value = 30.02298
